# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (rId1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value = 598
$wsExhibit.Range("F9").Value = 8827
$wsExhibit.Range("F14").Value = 117
$wsExhibit.Range("F18").Value = 279
$wsExhibit.Range("F20").Value = 235
$wsExhibit.Range("F21").Value = 1066

# Sheet: 全部类型 (rId4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 598
$wsAll.Range("F11").Value = 8827
$wsAll.Range("F16").Value = 117
$wsAll.Range("F20").Value = 279
$wsAll.Range("F22").Value = 235
$wsAll.Range("F23").Value = 1066
